# DowloadMauPhieuNhap.xlsx edit script
# - Update the "Ngay" (date) field
# - Replace the product name / unit on the first item row
# - Insert 4 more item rows (rows 14-17) with their own data,
#   pushing the totals / footer rows down accordingly
# - Update the final total (Cong tien) value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update date text in merged cell A6:F6
$ws.Range("A6").Value = "Ngày: 09/04/2023"

# 2) Insert four new rows right after row 13 (before the old "Cong tien" row).
#    This shifts the old rows 14-19 down to rows 18-23 and keeps their
#    formatting / merged ranges intact.
$ws.Rows("14:17").Insert()

# 3) Copy the formatting of row 13 onto the freshly inserted rows so the
#    new item rows look the same as the first item row.
$ws.Range("A13:H13").Copy()
$ws.Range("A14:H17").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows("14:17").RowHeight = 21

# 4) Update row 13 (first item) values
$ws.Range("B13").Value = "Tấm gương của những người thành đạt"
$ws.Range("C13").Value = 1810
$ws.Range("D13").Value = "Cái"
$ws.Range("E13").Value = 2000
$ws.Range("F13").Value = 3258000

# 5) Fill in the new item rows (14-17) - same product/unit text as row 13
$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "Tấm gương của những người thành đạt"
$ws.Range("C14").Value = 8
$ws.Range("D14").Value = "Cái"
$ws.Range("E14").Value = 2000
$ws.Range("F14").Value = 14400

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "Tấm gương của những người thành đạt"
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = "Cái"
$ws.Range("E15").Value = 2000
$ws.Range("F15").Value = 18000

$ws.Range("A16").Value = 4
$ws.Range("B16").Value = "Tấm gương của những người thành đạt"
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = "Cái"
$ws.Range("E16").Value = 2000
$ws.Range("F16").Value = 18000

$ws.Range("A17").Value = 5
$ws.Range("B17").Value = "Tấm gương của những người thành đạt"
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = "Cái"
$ws.Range("E17").Value = 2000
$ws.Range("F17").Value = 3600

# 6) Update the total amount (Cong tien row, now row 18)
$ws.Range("F18").Value = 3312000
